$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 (old duplicate "Contact") becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Remove the old duplicate "Contact" row (row 11), shifting rows 12-21 up
$meta.Rows("11").Delete()

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition updated to the specific value description
$elements.Range("K2").Value = "ReengagementValue"
$elements.Range("L2").Value = "Communication reengagement sequence value"
